$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.280225872993469
$ws.Range("B1").Value = 3.739739418029785
$ws.Range("C1").Value = 3.389606952667236
$ws.Range("D1").Value = 2.539255142211914
$ws.Range("E1").Value = 0.9886522889137268
